$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 111747705
$ws.Range("B2").Value = 93067
$ws.Range("D2").Value = 'LC'
$ws.Range("E2").Value = 2810
$ws.Range("F2").Value = 'Västlig hakmossa'
$ws.Range("G2").Value = 'Rhytidiadelphus loreus'
$ws.Range("H2").Value = '(Hedw.) Warnst.'
$ws.Range("P2").Value = 'Lilla gruvan (Lilla gruvan), Ög'
$ws.Range("Q2").Value = 575459
$ws.Range("R2").Value = 6509864
$ws.Range("S2").Value = 2
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()

# Row 3
$ws.Range("A3").Value = 111749006
$ws.Range("B3").Value = 8377
$ws.Range("D3").Value = 'LC'
$ws.Range("E3").Value = 106545
$ws.Range("F3").Value = 'Mindre märgborre'
$ws.Range("G3").Value = 'Tomicus minor'
$ws.Range("H3").Value = '(Hartig, 1834)'
$ws.Range("Q3").Value = 575512
$ws.Range("R3").Value = 6509826
$ws.Range("Z3").ClearContents()
$ws.Range("AB3").ClearContents()

# Row 4
$ws.Range("A4").Value = 111749860
$ws.Range("Q4").Value = 575357
$ws.Range("R4").Value = 6509772
$ws.Range("Z4").ClearContents()
$ws.Range("AB4").ClearContents()

# Row 5
$ws.Range("A5").Value = 111749343
$ws.Range("B5").Value = 78107
$ws.Range("D5").Value = 'NT'
$ws.Range("E5").Value = 6453
$ws.Range("F5").Value = 'Vedskivlav'
$ws.Range("G5").Value = 'Hertelidea botryosa'
$ws.Range("H5").Value = '(Fr.) Printzen & Kantvilas'
$ws.Range("Q5").Value = 575415
$ws.Range("R5").Value = 6509808
$ws.Range("S5").Value = 1
$ws.Range("Z5").ClearContents()
$ws.Range("AB5").ClearContents()

# Row 6
$ws.Range("A6").Value = 111749883
$ws.Range("B6").Value = 78107
$ws.Range("D6").Value = 'NT'
$ws.Range("E6").Value = 6453
$ws.Range("F6").Value = 'Vedskivlav'
$ws.Range("G6").Value = 'Hertelidea botryosa'
$ws.Range("H6").Value = '(Fr.) Printzen & Kantvilas'
$ws.Range("P6").Value = 'Älgsjöhåll (Älgsjöhåll), Ög'
$ws.Range("Q6").Value = 575337
$ws.Range("R6").Value = 6509789
$ws.Range("S6").Value = 1
$ws.Range("Z6").ClearContents()
$ws.Range("AB6").ClearContents()

# Row 7
$ws.Range("A7").Value = 111749897
$ws.Range("P7").Value = 'Älgsjöhåll (Älgsjöhåll), Ög'
$ws.Range("Q7").Value = 575337
$ws.Range("R7").Value = 6509781
$ws.Range("Z7").ClearContents()
$ws.Range("AB7").ClearContents()

# Row 8
$ws.Range("A8").Value = 111747186
$ws.Range("P8").Value = 'Lilla gruvan (Lilla gruvan), Ög'
$ws.Range("Q8").Value = 575436
$ws.Range("R8").Value = 6509857
$ws.Range("S8").Value = 2
$ws.Range("Z8").ClearContents()
$ws.Range("AB8").ClearContents()

# Row 9
$ws.Range("A9").Value = 111749097
$ws.Range("B9").Value = 93388
$ws.Range("E9").Value = 2180
$ws.Range("F9").Value = 'Blåmossa'
$ws.Range("G9").Value = 'Leucobryum glaucum'
$ws.Range("H9").Value = '(Hedw.) Ångstr.'
$ws.Range("Q9").Value = 575502
$ws.Range("R9").Value = 6509776
$ws.Range("S9").Value = 3
$ws.Range("Z9").ClearContents()
$ws.Range("AB9").ClearContents()

# Row 10
$ws.Range("A10").Value = 111964494
$ws.Range("B10").Value = 56414
$ws.Range("D10").Value = 'NT'
$ws.Range("E10").Value = 100049
$ws.Range("F10").Value = 'Spillkråka'
$ws.Range("G10").Value = 'Dryocopus martius'
$ws.Range("H10").Value = '(Linnaeus, 1758)'
$ws.Range("I10").ClearContents()
$ws.Range("J10").ClearContents()
$ws.Range("K10").Value = 'adult'
$ws.Range("M10").Value = 'förbiflygande'
$ws.Range("Q10").Value = 575346
$ws.Range("R10").Value = 6509958
$ws.Range("Z10").Value = '10:30'
$ws.Range("AB10").Value = '10:30'

# Row 11
$ws.Range("A11").Value = 111964550
$ws.Range("B11").Value = 103288
$ws.Range("D11").Value = 'LC'
$ws.Range("E11").Value = 221144
$ws.Range("F11").Value = 'Grönpyrola'
$ws.Range("G11").Value = 'Pyrola chlorantha'
$ws.Range("H11").Value = 'Sw.'
$ws.Range("I11").Value = "'25"
$ws.Range("J11").Value = 'plantor/tuvor'
$ws.Range("K11").Value = 'överblommad'
$ws.Range("M11").ClearContents()
$ws.Range("Q11").Value = 575346
$ws.Range("R11").Value = 6509958
$ws.Range("Z11").ClearContents()
$ws.Range("AB11").ClearContents()

# Row 12
$ws.Range("Q12").Value = 575609
$ws.Range("R12").Value = 6509825
$ws.Range("Z12").ClearContents()
$ws.Range("AB12").ClearContents()
